# "dashboard and course updated"
# Add a new data row (row 4) to Sheet1, mirroring the structure/styles of
# the existing rows 2-3, with a new "blog" card string in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: date value, formatted like A2/A3 (numFmtId 15 date style) ---
$ws.Range("A4").Value = 43895
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat

# --- Columns B,C,D,E,F,H,I: reuse text already present elsewhere on the sheet ---
# (use Value2 for reading -- Value returns a COM property descriptor in this host)
$ws.Range("B4").Value = $ws.Range("B3").Value2   # type: person
$ws.Range("C4").Value = $ws.Range("C3").Value2   # type: blog ser:62
$ws.Range("D4").Value = $ws.Range("D3").Value2   # type: course
$ws.Range("E4").Value = $ws.Range("E3").Value2   # type: meetup
$ws.Range("F4").Value = $ws.Range("F2").Value2   # type: signin
$ws.Range("H4").Value = $ws.Range("H2").Value2   # type: subscribe
$ws.Range("I4").Value = $ws.Range("I2").Value2   # type: footer

# --- Column G: brand-new shared string value ---
$ws.Range("G4").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 60"

# --- Columns J and K exist as empty styled cells on row 4 ---
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""

# --- Apply the same wrap-text style used on the rest of row 2/3 (B:K) ---
$ws.Range("B4:K4").WrapText = $true

# --- Row height matches the other data rows (409.6 = Excel's max row height) ---
$ws.Rows.Item(4).RowHeight = 409.6

# --- Update sheet view to reflect the newly-added/selected row, as in the diff ---
# (topLeftCell -> A4, i.e. row 4 / column A)
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G4").Select()
